$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Row 5: add new bug entry B-4 with description in columns E/F (F5 uses the
# same wrap-text/border style as the other description cells, e.g. F4)
$ws.Range("E5").Value = "B-4"
$ws.Range("F5").Value = "Бот игнорирует ввод отрицательного значения валюты и конвертирует его"
$ws.Range("F4").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 7: mark as "Not pass" and link to bug B-4
$ws.Range("B7").Value = "Not pass"
$ws.Range("C7").Value = "B-4"

# Row 5 now wraps two lines of text like row 4, so match its row height
$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(4).RowHeight

# Update selection to F8 as in the diff
$ws.Range("F8").Select()
